# ---------------------------------------------------------------------
# "Add files via upload" commit — the document's text is reworded so
# that the generic "ενδιαφερόμενος" ("interested party") becomes the
# more specific "ιδιώτης" ("private individual") everywhere it is used
# as the actor in the use-case steps. Every grammatical form that
# literally spells "ενδιαφερόμενος" (nominative) is replaced; the
# genitive/accusative forms ("ενδιαφερόμενου", "ενδιαφερόμενο") are
# distinct words and are intentionally left untouched.
# ---------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. Cosmetic split of the title run (no visible text change) ------
# In the target file "Use Cases" is split into "Use" + " " + "Cases"
# (the middle space loses the run's language tag). Reproduce the run
# split with a harmless formatting no-op (Bold on/off) on the space
# character so Word carves it into its own run; this does not alter
# the rendered text.
$title = $d.Paragraphs(1).Range
if ($title.Text.TrimEnd([char]13) -eq "Use Cases") {
    $spacePos = $title.Start + 3
    $rSpace = $d.Range($spacePos, $spacePos + 1)
    if ($rSpace.Text -eq " ") {
        $rSpace.Font.Bold = $true
        $rSpace.Font.Bold = $false
    }
}

# --- 2. Reword "ενδιαφερόμενος" -> "ιδιώτης" everywhere ----------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "ενδιαφερόμενος",   # FindText
    $false,              # MatchCase
    $false,              # MatchWholeWord
    $false,              # MatchWildcards
    $false,              # MatchSoundsLike
    $false,              # MatchAllWordForms
    $true,               # Forward
    1,                   # Wrap (wdFindContinue)
    $false,              # Format
    "ιδιώτης",           # ReplaceWith
    2                    # Replace (wdReplaceAll)
)

Write-Output "ok"
